$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "group" column (bold, matching the other header cells)
$ws.Range("D1").Value = "group"
$ws.Range("D1").Font.Bold = $true

# Fill column D for existing data rows (2-60) with "None"
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 4).Value = "None"
}

# Row 61 (bean-action-potentials) gets a special group value
$ws.Cells.Item(61, 4).Value = "NURB 31800"

# Row 62 (natural-image-statistics) gets "None"
$ws.Cells.Item(62, 4).Value = "None"

# New row 63: Elements of Information Theory
$ws.Cells.Item(63, 1).Value = "information-theory-elements"
$ws.Cells.Item(63, 2).Value = "Elements of Information Theory"
$ws.Cells.Item(63, 3).Value = "information-theory"
$ws.Cells.Item(63, 4).Value = "None"

# Update the window scroll position / active selection to match the
# author's final cursor position on the sheet.
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D61").Select()
